{"js": "// Split the four concatenated \"Programa resumido\"/\"Programa\" summary\n// paragraphs (PT + EN, short + long) into line-broken runs: each\n// numbered item becomes its own <w:t> segment separated by a manual\n// line break (<w:br/>), matching the authored OOXML diff.\nconst EDITS = [\n  {\n    original: \"1 \u2013 Introdu\u00e7\u00e3o aos Sistemas Produtivos2 \u2013 Papel Estrat\u00e9gico da Produ\u00e7\u00e3o3 \u2013 Estrat\u00e9gia de Produ\u00e7\u00e3o4 \u2013 Projeto em Gest\u00e3o de Produ\u00e7\u00e3o5 \u2013 Projeto de Produtos e Servi\u00e7os6 \u2013 Projeto da Rede de Opera\u00e7\u00f5es Produtivas7 \u2013 Arranjo F\u00edsico e Fluxo\",\n    segments: [\n      \"1 \u2013 Introdu\u00e7\u00e3o aos Sistemas Produtivos\",\n      \"2 \u2013 Papel Estrat\u00e9gico da Produ\u00e7\u00e3o\",\n      \"3 \u2013 Estrat\u00e9gia de Produ\u00e7\u00e3o\",\n      \"4 \u2013 Projeto em Gest\u00e3o de Produ\u00e7\u00e3o\",\n      \"5 \u2013 Projeto de Produtos e Servi\u00e7os\",\n      \"6 \u2013 Projeto da Rede de Opera\u00e7\u00f5es Produtivas\",\n      \"7 \u2013 Arranjo F\u00edsico e Fluxo\",\n    ],\n  },\n  {\n    original: \"1 - Introduction to Productive Systems2 - Strategic Role of Production3 - Production Strategy4 - Project in Production Management5 - Product and Service Project6 - Production Operations Network Project7 - Physical Arrangement and Flow\",\n    segments: [\n      \"1 - Introduction to Productive Systems\",\n      \"2 - Strategic Role of Production\",\n      \"3 - Production Strategy\",\n      \"4 - Project in Production Management\",\n      \"5 - Product and Service Project\",\n      \"6 - Production Operations Network Project\",\n      \"7 - Physical Arrangement and Flow\",\n    ],\n  },\n  {\n    original: \"1 \u2013 Introdu\u00e7\u00e3o aos Sistemas Produtivos; Produ\u00e7\u00e3o na Organiza\u00e7\u00e3o. Inputs, Processos de Transforma\u00e7\u00e3o e Outputs. Tipos de Opera\u00e7\u00f5es de Produ\u00e7\u00e3o. Atividades da administra\u00e7\u00e3o da produ\u00e7\u00e3o.2 \u2013 Papel Estrat\u00e9gico da Produ\u00e7\u00e3o; Papel da fun\u00e7\u00e3o produ\u00e7\u00e3o. Objetivos de desempenho.3 \u2013 Tipos de Manufatura; Tipos b\u00e1sicos de Manufatura.4 \u2013 Arranjo F\u00edsico e Fluxo; Procedimento de Arranjo F\u00edsico. Tipos b\u00e1sicos de arranjo f\u00edsico. Projeto de arranjo f\u00edsico.5 \u2013 Organiza\u00e7\u00e3o do Trabalho e M\u00e9todos;T\u00e9cnicas de organiza\u00e7\u00e3o e m\u00e9todos de trabalho6 - Introdu\u00e7\u00e3o ao Planejamento e Controle de Produ\u00e7\u00e3o.Conceitua\u00e7\u00e3o do PCP; concilia\u00e7\u00e3o de suprimento e demanda; natureza do suprimento e da demanda; atividades de PCP; efeito volume-variedade no PCP.7 - Introdu\u00e7\u00e3o \u00e0 qualidade e a tecnologia de processo Import\u00e2ncia; vis\u00f5es; princ\u00edpios de administra\u00e7\u00e3o da qualidade total.\",\n    segments: [\n      \"1 \u2013 Introdu\u00e7\u00e3o aos Sistemas Produtivos; Produ\u00e7\u00e3o na Organiza\u00e7\u00e3o. Inputs, Processos de Transforma\u00e7\u00e3o e Outputs. Tipos de Opera\u00e7\u00f5es de Produ\u00e7\u00e3o. Atividades da administra\u00e7\u00e3o da produ\u00e7\u00e3o.\",\n      \"2 \u2013 Papel Estrat\u00e9gico da Produ\u00e7\u00e3o; Papel da fun\u00e7\u00e3o produ\u00e7\u00e3o. Objetivos de desempenho.\",\n      \"3 \u2013 Tipos de Manufatura; Tipos b\u00e1sicos de Manufatura.\",\n      \"4 \u2013 Arranjo F\u00edsico e Fluxo; Procedimento de Arranjo F\u00edsico. Tipos b\u00e1sicos de arranjo f\u00edsico. Projeto de arranjo f\u00edsico.\",\n      \"5 \u2013 Organiza\u00e7\u00e3o do Trabalho e M\u00e9todos;T\u00e9cnicas de organiza\u00e7\u00e3o e m\u00e9todos de trabalho\",\n      \"6 - Introdu\u00e7\u00e3o ao Planejamento e Controle de Produ\u00e7\u00e3o.\",\n      \"Conceitua\u00e7\u00e3o do PCP; concilia\u00e7\u00e3o de suprimento e demanda; natureza do suprimento e da demanda; atividades de PCP; efeito volume-variedade no PCP.\",\n      \"7 - Introdu\u00e7\u00e3o \u00e0 qualidade e a tecnologia de processo Import\u00e2ncia; vis\u00f5es; princ\u00edpios de administra\u00e7\u00e3o da qualidade total.\",\n    ],\n  },\n  {\n    original: \"1 \u2013 Introduction to Production Systems;Production in the Organization. Inputs, Transformation Processes and Outputs. Types of Production Operations. Production management activities.2 \u2013 Strategic Role of Production;Role of the production function. Performance objectives.3 \u2013 Types of Manufacturing;Basic types of Manufacturing.4 \u2013 Physical Arrangement and Flow;Physical Arrangement Procedure. Basic types of physical arrangement. Physical arrangement design.5 \u2013 Work Organization and Methods;Organization techniques and work methods6 - Introduction to Production Planning and Control.Conceptualization of the PCP; reconciliation of supply and demand; nature of supply and demand; PCP activities; volume-variety effect in PCP.7 - Introduction to quality and process technologyImportance; visions; total quality management principles.\",\n    segments: [\n      \"1 \u2013 Introduction to Production Systems;\",\n      \"Production in the Organization. Inputs, Transformation Processes and Outputs. Types of Production Operations. Production management activities.\",\n      \"2 \u2013 Strategic Role of Production;\",\n      \"Role of the production function. Performance objectives.\",\n      \"3 \u2013 Types of Manufacturing;\",\n      \"Basic types of Manufacturing.\",\n      \"4 \u2013 Physical Arrangement and Flow;\",\n      \"Physical Arrangement Procedure. Basic types of physical arrangement. Physical arrangement design.\",\n      \"5 \u2013 Work Organization and Methods;\",\n      \"Organization techniques and work methods\",\n      \"6 - Introduction to Production Planning and Control.\",\n      \"Conceptualization of the PCP; reconciliation of supply and demand; nature of supply and demand; PCP activities; volume-variety effect in PCP.\",\n      \"7 - Introduction to quality and process technology\",\n      \"Importance; visions; total quality management principles.\",\n    ],\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nfor (const edit of EDITS) {\n  const target = paragraphs.items.find((p) => p.text === edit.original);\n  if (!target) {\n    throw new Error(\"Paragraph not found for: \" + edit.original.slice(0, 60));\n  }\n  // U+000B (vertical tab) is how the Word object model represents a\n  // manual line break inside Range/Paragraph text; inserting it via\n  // insertText keeps everything in one run and emits <w:t>..</w:t><w:br/>\n  // pairs in place of the break, exactly like the target markup.\n  const newText = edit.segments.join(\"\\u000b\");\n  target.getRange().insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Split the four concatenated \"Programa resumido\"/\"Programa\" summary\n# paragraphs (PT + EN, short + long) into line-broken runs: each\n# numbered item becomes its own text segment separated by a manual\n# line break (vertical tab / Chr(11), which Word renders as <w:br/>),\n# matching the authored OOXML diff.\n\n$vt = [char]11  # manual line break char used by Word's Range.Text / Selection.TypeText\n\n$edits = @(\n    [PSCustomObject]@{\n        Original = '1 \u2013 Introdu\u00e7\u00e3o aos Sistemas Produtivos2 \u2013 Papel Estrat\u00e9gico da Produ\u00e7\u00e3o3 \u2013 Estrat\u00e9gia de Produ\u00e7\u00e3o4 \u2013 Projeto em Gest\u00e3o de Produ\u00e7\u00e3o5 \u2013 Projeto de Produtos e Servi\u00e7os6 \u2013 Projeto da Rede de Opera\u00e7\u00f5es Produtivas7 \u2013 Arranjo F\u00edsico e Fluxo'\n        Segments = @(\n            '1 \u2013 Introdu\u00e7\u00e3o aos Sistemas Produtivos',\n            '2 \u2013 Papel Estrat\u00e9gico da Produ\u00e7\u00e3o',\n            '3 \u2013 Estrat\u00e9gia de Produ\u00e7\u00e3o',\n            '4 \u2013 Projeto em Gest\u00e3o de Produ\u00e7\u00e3o',\n            '5 \u2013 Projeto de Produtos e Servi\u00e7os',\n            '6 \u2013 Projeto da Rede de Opera\u00e7\u00f5es Produtivas',\n            '7 \u2013 Arranjo F\u00edsico e Fluxo'\n        )\n    },\n    [PSCustomObject]@{\n        Original = '1 - Introduction to Productive Systems2 - Strategic Role of Production3 - Production Strategy4 - Project in Production Management5 - Product and Service Project6 - Production Operations Network Project7 - Physical Arrangement and Flow'\n        Segments = @(\n            '1 - Introduction to Productive Systems',\n            '2 - Strategic Role of Production',\n            '3 - Production Strategy',\n            '4 - Project in Production Management',\n            '5 - Product and Service Project',\n            '6 - Production Operations Network Project',\n            '7 - Physical Arrangement and Flow'\n        )\n    },\n    [PSCustomObject]@{\n        Original = '1 \u2013 Introdu\u00e7\u00e3o aos Sistemas Produtivos; Produ\u00e7\u00e3o na Organiza\u00e7\u00e3o. Inputs, Processos de Transforma\u00e7\u00e3o e Outputs. Tipos de Opera\u00e7\u00f5es de Produ\u00e7\u00e3o. Atividades da administra\u00e7\u00e3o da produ\u00e7\u00e3o.2 \u2013 Papel Estrat\u00e9gico da Produ\u00e7\u00e3o; Papel da fun\u00e7\u00e3o produ\u00e7\u00e3o. Objetivos de desempenho.3 \u2013 Tipos de Manufatura; Tipos b\u00e1sicos de Manufatura.4 \u2013 Arranjo F\u00edsico e Fluxo; Procedimento de Arranjo F\u00edsico. Tipos b\u00e1sicos de arranjo f\u00edsico. Projeto de arranjo f\u00edsico.5 \u2013 Organiza\u00e7\u00e3o do Trabalho e M\u00e9todos;T\u00e9cnicas de organiza\u00e7\u00e3o e m\u00e9todos de trabalho6 - Introdu\u00e7\u00e3o ao Planejamento e Controle de Produ\u00e7\u00e3o.Conceitua\u00e7\u00e3o do PCP; concilia\u00e7\u00e3o de suprimento e demanda; natureza do suprimento e da demanda; atividades de PCP; efeito volume-variedade no PCP.7 - Introdu\u00e7\u00e3o \u00e0 qualidade e a tecnologia de processo Import\u00e2ncia; vis\u00f5es; princ\u00edpios de administra\u00e7\u00e3o da qualidade total.'\n        Segments = @(\n            '1 \u2013 Introdu\u00e7\u00e3o aos Sistemas Produtivos; Produ\u00e7\u00e3o na Organiza\u00e7\u00e3o. Inputs, Processos de Transforma\u00e7\u00e3o e Outputs. Tipos de Opera\u00e7\u00f5es de Produ\u00e7\u00e3o. Atividades da administra\u00e7\u00e3o da produ\u00e7\u00e3o.',\n            '2 \u2013 Papel Estrat\u00e9gico da Produ\u00e7\u00e3o; Papel da fun\u00e7\u00e3o produ\u00e7\u00e3o. Objetivos de desempenho.',\n            '3 \u2013 Tipos de Manufatura; Tipos b\u00e1sicos de Manufatura.',\n            '4 \u2013 Arranjo F\u00edsico e Fluxo; Procedimento de Arranjo F\u00edsico. Tipos b\u00e1sicos de arranjo f\u00edsico. Projeto de arranjo f\u00edsico.',\n            '5 \u2013 Organiza\u00e7\u00e3o do Trabalho e M\u00e9todos;T\u00e9cnicas de organiza\u00e7\u00e3o e m\u00e9todos de trabalho',\n            '6 - Introdu\u00e7\u00e3o ao Planejamento e Controle de Produ\u00e7\u00e3o.',\n            'Conceitua\u00e7\u00e3o do PCP; concilia\u00e7\u00e3o de suprimento e demanda; natureza do suprimento e da demanda; atividades de PCP; efeito volume-variedade no PCP.',\n            '7 - Introdu\u00e7\u00e3o \u00e0 qualidade e a tecnologia de processo Import\u00e2ncia; vis\u00f5es; princ\u00edpios de administra\u00e7\u00e3o da qualidade total.'\n        )\n    },\n    [PSCustomObject]@{\n        Original = '1 \u2013 Introduction to Production Systems;Production in the Organization. Inputs, Transformation Processes and Outputs. Types of Production Operations. Production management activities.2 \u2013 Strategic Role of Production;Role of the production function. Performance objectives.3 \u2013 Types of Manufacturing;Basic types of Manufacturing.4 \u2013 Physical Arrangement and Flow;Physical Arrangement Procedure. Basic types of physical arrangement. Physical arrangement design.5 \u2013 Work Organization and Methods;Organization techniques and work methods6 - Introduction to Production Planning and Control.Conceptualization of the PCP; reconciliation of supply and demand; nature of supply and demand; PCP activities; volume-variety effect in PCP.7 - Introduction to quality and process technologyImportance; visions; total quality management principles.'\n        Segments = @(\n            '1 \u2013 Introduction to Production Systems;',\n            'Production in the Organization. Inputs, Transformation Processes and Outputs. Types of Production Operations. Production management activities.',\n            '2 \u2013 Strategic Role of Production;',\n            'Role of the production function. Performance objectives.',\n            '3 \u2013 Types of Manufacturing;',\n            'Basic types of Manufacturing.',\n            '4 \u2013 Physical Arrangement and Flow;',\n            'Physical Arrangement Procedure. Basic types of physical arrangement. Physical arrangement design.',\n            '5 \u2013 Work Organization and Methods;',\n            'Organization techniques and work methods',\n            '6 - Introduction to Production Planning and Control.',\n            'Conceptualization of the PCP; reconciliation of supply and demand; nature of supply and demand; PCP activities; volume-variety effect in PCP.',\n            '7 - Introduction to quality and process technology',\n            'Importance; visions; total quality management principles.'\n        )\n    }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($edit in $edits) {\n    $found = $false\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs($i)\n        # Paragraph.Range.Text carries a trailing paragraph mark (Chr(13)),\n        # so trim it before comparing against the diff's plain run text.\n        $text = $p.Range.Text.TrimEnd([char]13)\n        if ($text -eq $edit.Original) {\n            $p.Range.Text = [string]::Join($vt, $edit.Segments)\n            $found = $true\n            break\n        }\n    }\n    if (-not $found) {\n        throw \"Paragraph not found for: \" + $edit.Original.Substring(0, 60)\n    }\n}\n"}
